$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 2
$ws.Range("A2").Value = 800
$ws.Range("B2").Value = 900

# Reset the view so the top-left visible cell is A1 and B3 is selected
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B3").Select()
